$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$ws.Range("C2").Value = 288
$ws.Range("D2").Value = 228.3

$ws.Range("C3").Value = 352
$ws.Range("D3").Value = 215.96

$ws.Range("C4").Value = 190
$ws.Range("D4").Value = 212.04

$ws.Range("C5").Value = 207
$ws.Range("D5").Value = 239.79

$ws.Range("C6").Value = 227
$ws.Range("D6").Value = 242.28

$ws.Range("C7").Value = 111
$ws.Range("D7").Value = 114.15

$ws.Range("C8").Value = 86
$ws.Range("D8").Value = 92.46
